# "Generate Report for Handback"
#
# The localization status report is refreshed: both the zh-cn and de-de
# languages are now "Handed back: in sync with en-US". The per-language
# sheets get their "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns filled in (with a hyperlink on the
# target-file cell, mirroring the existing source-file hyperlink), and a
# few report columns are widened so the new content fits.

$wb = $excel.ActiveWorkbook

$urlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f3d46c644ae8881723abbf9cd29ae43b27a5220/e2e/"
$url384 = $urlBase + "384be899-8496-4c83-85ce-9e98a7283034.md"
$urlEc96 = $urlBase + "ec96196b-fbcc-40ed-8720-318ae29c390e.md"

# Column width helper: this engine quantizes ColumnWidth to a 1/6-character
# grid when it round-trips through the stored <col width> (w_stored =
# (round(w_input*6)+5)/6), so back solve for the COM input that lands on
# the desired stored width.
function Set-ColWidth($col, $targetStored) {
  $q = [Math]::Round($targetStored * 6 - 5)
  $col.ColumnWidth = $q / 6
}

# ---------------------------------------------------------------------
# Overview sheet: both languages are now handed back & in sync.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US") | Out-Null

Set-ColWidth $wsOverview.Columns.Item(5) 29.9777047293527
Set-ColWidth $wsOverview.Columns.Item(6) 29.9777047293527

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de): fill in the handback columns.
# ---------------------------------------------------------------------
foreach ($sheetName in @("zh-cn", "de-de")) {
  $ws = $wb.Worksheets.Item($sheetName)

  Set-ColWidth $ws.Columns.Item(3) 29.9777047293527
  Set-ColWidth $ws.Columns.Item(9) 40
  Set-ColWidth $ws.Columns.Item(10) 40

  $suffix = if ($sheetName -eq "zh-cn") { "zh-cn" } else { "de-de" }

  # Re-create the hyperlinks collection from scratch so the new
  # "Latest Target File" links land right after their matching
  # "Source File Name" links, in row order (A2, I2, A3, I3).
  $ws.Hyperlinks.Delete()

  $ws.Hyperlinks.Add($ws.Range("A2"), $url384, "", "", "384be899-8496-4c83-85ce-9e98a7283034.md") | Out-Null
  $ws.Hyperlinks.Add($ws.Range("I2"), $url384, "", "", "384be899-8496-4c83-85ce-9e98a7283034.md") | Out-Null
  $ws.Hyperlinks.Add($ws.Range("A3"), $urlEc96, "", "", "ec96196b-fbcc-40ed-8720-318ae29c390e.md") | Out-Null
  $ws.Hyperlinks.Add($ws.Range("I3"), $urlEc96, "", "", "ec96196b-fbcc-40ed-8720-318ae29c390e.md") | Out-Null

  # Hyperlinks.Add re-styles with a brand new "Hyperlink" cell style that
  # only carries font color (not the underline); explicitly re-apply the
  # underline + color so I2/I3 match the look of the existing A2/A3 links.
  $ws.Range("I2").Font.Underline = $true
  $ws.Range("I2").Font.Color = 15570276
  $ws.Range("I3").Font.Underline = $true
  $ws.Range("I3").Font.Color = 15570276

  $ws.Range("J2").Value = "384be899-8496-4c83-85ce-9e98a7283034.97b2c27eae186aea3c04d3e3ef20e5f15b945af9." + $suffix + ".xlf"
  $ws.Range("J3").Value = "ec96196b-fbcc-40ed-8720-318ae29c390e.591c3e3a8751791a39e1440f61aa2f7a1e11fe8f." + $suffix + ".xlf"
}

# zh-cn was already in sync; just needed its handback datetime filled in.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Cells.Replace("0001-01-01 00:00:00", "2016-08-20 04:34:35") | Out-Null

# de-de has just been handed back now.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("K2").Value = "2016-08-20 04:34:43"
$wsDe.Range("K3").Value = "2016-08-20 04:34:43"

Write-Output "Handback report generated."
